$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.700.60'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.895.87'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4922'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.69%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2938'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06744'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.896.50'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '17.22'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07247'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.69%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '90.86'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6760'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.92%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.031'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.696.86'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007974'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.66%  '
$ws.Range('E18').Value = '  +0.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.09'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.141.34'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.808'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '191.64'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +34.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.084'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.367'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '156.39'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +12.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.894'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.403'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.297'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09073'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.56%  '
$ws.Range('E32').Value = '  +0.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05218'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7398'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.107'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('E36').Value = '  +3.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01832'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.679'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.124'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9310'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4393'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '105.08'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.40%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.723'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.536'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.16%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1351'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05862'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.740'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.37%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.3933'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.19%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.66'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.74%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.416'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.60%  '
